$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 'Price' column (D) holds plain numeric-looking text, e.g. "515.51".
# Assigning such a string straight to .Value lets Excel auto-convert it to a
# real number, which would change the cell's stored type. Temporarily mark the
# whole column as Text before writing the new values, then restore the original
# (default) style afterwards so no formatting differences are left behind.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '56.253.30'
$ws.Range('E2').Value = '  +2.12%  '
$ws.Range('D3').Value = '2.319.79'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').Value = '  +0.60%  '
$ws.Range('D5').Value = '515.51'
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('D6').Value = '133.30'
$ws.Range('E6').Value = '  +3.26%  '
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('D8').Value = '0.534'
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').Value = '2.341.28'
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('E10').Value = '  +5.89%  '
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').Value = '5.22'
$ws.Range('E12').Value = '  +5.80%  '
$ws.Range('D13').Value = '0.338'
$ws.Range('E13').Value = '  -2.24%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '23.62'
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.738.24'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '56.518.68'
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').Value = '2.332.99'
$ws.Range('E18').Value = '  +1.33%  '
$ws.Range('D19').Value = '10.37'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').Value = '4.24'
$ws.Range('E20').Value = '  +2.47%  '
$ws.Range('D21').Value = '319.26'
$ws.Range('E21').Value = '  +3.81%  '
$ws.Range('D22').Value = '6.62'
$ws.Range('E22').Value = '  +2.99%  '
$ws.Range('D23').Value = '0.996'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').Value = '60.44'
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').Value = '0.158'
$ws.Range('E26').Value = '  +4.67%  '
$ws.Range('E27').Value = '  +3.81%  '
$ws.Range('D28').Value = '170.90'
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('E29').Value = '  +9.34%  '
$ws.Range('D30').Value = '0.0₃0733'
$ws.Range('E30').Value = '  +4.00%  '
$ws.Range('D31').Value = '6.21'
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('E32').Value = '  +2.32%  '
$ws.Range('D33').Value = '18.22'
$ws.Range('E33').Value = '  +1.34%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = '0.993'
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('D36').Value = '0.947'
$ws.Range('E36').Value = '  +2.86%  '
$ws.Range('E37').Value = '  +2.94%  '
$ws.Range('D38').Value = '3.96'
$ws.Range('E38').Value = '  +4.60%  '
$ws.Range('D39').Value = '1.52'
$ws.Range('E39').Value = '  +7.37%  '
$ws.Range('D40').Value = '37.42'
$ws.Range('E40').Value = '  +2.81%  '
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('D42').Value = '137.81'
$ws.Range('E42').Value = '  +9.02%  '
$ws.Range('D43').Value = '3.55'
$ws.Range('E43').Value = '  +4.02%  '
$ws.Range('D44').Value = '275.01'
$ws.Range('E44').Value = '  +10.27%  '
$ws.Range('D45').Value = '5.02'
$ws.Range('E45').Value = '  +1.77%  '
$ws.Range('D46').Value = '0.0928'
$ws.Range('E46').Value = '  +2.57%  '
$ws.Range('D47').Value = '0.0503'
$ws.Range('E47').Value = '  +0.93%  '
$ws.Range('E48').Value = '  +1.33%  '
$ws.Range('D49').Value = '0.0215'
$ws.Range('E49').Value = '  +3.85%  '
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('D51').Value = '16.73'
$ws.Range('E51').Value = '  +1.77%  '

$ws.Range('D2:D51').Style = 'Normal'
